# Corrected CO2 and GHG Emissions
# Update the raw input values (co2 = G, ch4 = J, n2o = M, ghg_100 = AE,
# ghg_20 = AH columns) for the Fuel Cell, Gas Turbine, Microturbine and RE
# technology blocks, plus the "ABC Only" n2o subtotal row. All dependent
# formulas (percentages, sums, totals) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fuel Cell (rows 4-5) ---
$ws.Range("G4").Value = 1140
$ws.Range("J4").Value = 479
$ws.Range("AE4").Value = 1119
$ws.Range("AH4").Value = 736

$ws.Range("G5").Value = 542
$ws.Range("J5").Value = 200
$ws.Range("M5").Value = 829
$ws.Range("AE5").Value = 457
$ws.Range("AH5").Value = 263

# --- Gas Turbine (rows 6-7) ---
$ws.Range("G7").Value = 10
$ws.Range("M7").Value = 1187

# --- Microturbine (rows 8-9) ---
$ws.Range("G9").Value = 1
$ws.Range("M9").Value = 1270

# --- RE (rows 10-11) ---
$ws.Range("G10").Value = 1047
$ws.Range("AE10").Value = 666
$ws.Range("AH10").Value = 240

$ws.Range("G11").Value = 488
$ws.Range("J11").Value = 31
$ws.Range("M11").Value = 926
$ws.Range("AE11").Value = 258
$ws.Range("AH11").Value = 78

# --- ABC Only (row 12) ---
$ws.Range("M12").Value = 0
